$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear forecast cells for the earliest rows (2-6) in columns C and E
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Update recalculated forecast values for rows 7-19
$ws.Range("C7").Value = -0.6616365666142765
$ws.Range("E7").Value = 0.2197847717222867

$ws.Range("C8").Value = 0.2184978785563896
$ws.Range("E8").Value = 0.04195831742983547

$ws.Range("C9").Value = -0.01790997771649039
$ws.Range("E9").Value = 0.1671491311400208

$ws.Range("C10").Value = -0.08711135105702317
$ws.Range("E10").Value = 0.1158714888162216

$ws.Range("C11").Value = 0.1761917659537371
$ws.Range("E11").Value = 0.1078587431702305

$ws.Range("C12").Value = 0.3524405906205841
$ws.Range("E12").Value = 0.1531204771924033

$ws.Range("C13").Value = 0.07560805834034845
$ws.Range("E13").Value = -0.01670081902098719

$ws.Range("C14").Value = -0.5849047489490333
$ws.Range("E14").Value = -0.2059746096811033

$ws.Range("C15").Value = -0.4334047671505248
$ws.Range("E15").Value = -0.06335028919957075

$ws.Range("C16").Value = -0.1663214453978101
$ws.Range("E16").Value = -0.3613518455741316

$ws.Range("C17").Value = 0.6502606143725664
$ws.Range("E17").Value = -0.1691853834640433

$ws.Range("C18").Value = 0.1549171986535924
$ws.Range("E18").Value = -0.03486668218654065

$ws.Range("C19").Value = -0.06391119588061711
$ws.Range("E19").Value = -0.2122873162357264
